# Apply updated values to column E (row data) on sheet "GSI"
# Commit: correctif probleme insertion matiere dans bd et mise a jour semestre etudiant lors de l'inscription

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSI")

$updates = @{3 = 16; 4 = 12; 5 = 14; 6 = 8; 8 = 10; 9 = 5; 11 = 18; 12 = 20; 13 = 15; 14 = 19; 15 = 20; 16 = 17; 17 = 18; 18 = 18; 19 = 18; 20 = 5; 21 = 14; 22 = 16; 23 = 18; 24 = 19; 25 = 15; 26 = 16; 27 = 8; 28 = 16; 31 = 10; 32 = 18; 33 = 17; 34 = 9; 36 = 6; 37 = 8; 38 = 17; 39 = 16; 40 = 8; 41 = 9; 42 = 12; 43 = 16; 44 = 17; 45 = 12; 46 = 9; 47 = 10; 48 = 14; 49 = 12; 50 = 7; 51 = 18; 52 = 11; 53 = 8; 54 = 15; 55 = 19; 56 = 16; 57 = 10; 58 = 11; 59 = 8; 60 = 20; 61 = 16; 62 = 12; 63 = 16}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}
